# Apply "repull data, push all data, mean calculation" edits to the dSF (column F) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F ("dSF")
$updates = @{
    8  = 2
    9  = 5
    14 = -3
    22 = -2
    26 = 0
    28 = 1
    29 = 0
    34 = 2
    36 = 0
    41 = 1
    58 = -3
    59 = -1
    60 = -11
    64 = -1
    72 = 2
    74 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
